$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1980
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 1980
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 1980
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -2948

$ws.Range("H132").Value = 2977400.2
$ws.Range("I132").Value = 3126000.2
$ws.Range("J132").Value = 5400
$ws.Range("K132").Value = 9378000.600000001
$ws.Range("L132").Value = 16200
$ws.Range("M132").Value = -9375470.600000001
$ws.Range("N132").Value = -21260

$ws.Range("H135").Value = 1480.4412
$ws.Range("I135").Value = 547.38464
$ws.Range("J135").Value = 4512.875
$ws.Range("K135").Value = 4926.46176
$ws.Range("L135").Value = 40615.875
$ws.Range("M135").Value = -2391.46176
$ws.Range("N135").Value = -45685.875

$ws.Range("H138").Value = 2784.5374
$ws.Range("I138").Value = 1340.973
$ws.Range("J138").Value = 4564.933
$ws.Range("K138").Value = 4022.919
$ws.Range("L138").Value = 13694.799
$ws.Range("M138").Value = 1117.081
$ws.Range("N138").Value = -23974.799

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3004.988
$ws.Range("I32").Value = 2241.3333
$ws.Range("J32").Value = 10164.25
$ws.Range("K32").Value = 2241.3333
$ws.Range("L32").Value = 10164.25
$ws.Range("M32").Value = -1954.3333
$ws.Range("N32").Value = -10738.25

$ws.Range("H45").Value = 1257.0714
$ws.Range("I45").Value = 1257.0714
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1257.0714
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -880.0714
$ws.Range("N45").ClearContents()

$ws.Range("H61").Value = 1027.909
$ws.Range("I61").Value = 997.09375
$ws.Range("K61").Value = 997.09375
$ws.Range("M61").Value = -785.09375

$ws.Range("H74").Value = 1306.6562
$ws.Range("I74").Value = 865.4167
$ws.Range("J74").Value = 2630.375
$ws.Range("K74").Value = 865.4167
$ws.Range("L74").Value = 2630.375
$ws.Range("M74").Value = 8.583300000000008
$ws.Range("N74").Value = -4378.375

$ws.Range("H77").Value = 1306.6562
$ws.Range("I77").Value = 865.4167
$ws.Range("J77").Value = 2630.375
$ws.Range("K77").Value = 4327.0835
$ws.Range("L77").Value = 13151.875
$ws.Range("M77").Value = 40.91650000000027
$ws.Range("N77").Value = -21887.875

$ws.Range("H92").Value = 25968.092
$ws.Range("J92").Value = 25968.092
$ws.Range("L92").Value = 25968.092
$ws.Range("N92").Value = -30960.092

$ws.Range("H102").Value = 3850.36
$ws.Range("I102").Value = 4016.1875
$ws.Range("J102").Value = 3555.5557
$ws.Range("K102").Value = 4016.1875
$ws.Range("L102").Value = 3555.5557
$ws.Range("M102").Value = -2394.1875
$ws.Range("N102").Value = -6799.5557

$ws.Range("H132").Value = 2007.8096
$ws.Range("I132").Value = 1867.4359
$ws.Range("J132").Value = 3832.6667
$ws.Range("K132").Value = 5602.307699999999
$ws.Range("L132").Value = 11498.0001
$ws.Range("M132").Value = -3072.307699999999
$ws.Range("N132").Value = -16558.0001

$ws.Range("H133").Value = 50261
$ws.Range("J133").Value = 50261
$ws.Range("L133").Value = 50261
$ws.Range("N133").Value = -55321

$ws.Range("H136").Value = 1027.909
$ws.Range("I136").Value = 997.09375
$ws.Range("K136").Value = 2991.28125
$ws.Range("M136").Value = -441.28125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 44468
$ws.Range("J51").Value = 44468
$ws.Range("L51").Value = 44468
$ws.Range("N51").Value = -45450

$ws.Range("H94").Value = 685.42426
$ws.Range("I94").Value = 521.11536
$ws.Range("J94").Value = 1295.7142
$ws.Range("K94").Value = 521.11536
$ws.Range("L94").Value = 1295.7142
$ws.Range("M94").Value = -70.11536000000001
$ws.Range("N94").Value = -2197.7142

$ws.Range("H107").Value = 1301.5
$ws.Range("I107").Value = 961.1
$ws.Range("J107").Value = 1868.8334
$ws.Range("K107").Value = 961.1
$ws.Range("L107").Value = 1868.8334
$ws.Range("M107").Value = 958.9
$ws.Range("N107").Value = -5708.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20582.465
$ws.Range("I31").Value = 30011.555
$ws.Range("J31").Value = 5153.0454
$ws.Range("K31").Value = 30011.555
$ws.Range("L31").Value = 5153.0454
$ws.Range("M31").Value = -29716.555
$ws.Range("N31").Value = -5743.0454

$ws.Range("H34").Value = 20582.465
$ws.Range("I34").Value = 30011.555
$ws.Range("J34").Value = 5153.0454
$ws.Range("K34").Value = 30011.555
$ws.Range("L34").Value = 5153.0454
$ws.Range("M34").Value = -29809.555
$ws.Range("N34").Value = -5557.0454

$ws.Range("H132").Value = 1522.8462
$ws.Range("I132").Value = 890.8182
$ws.Range("K132").Value = 2672.4546
$ws.Range("M132").Value = -142.4546

$ws.Range("H134").Value = 1391.6207
$ws.Range("I134").Value = 1112.1482
$ws.Range("J134").Value = 5164.5
$ws.Range("K134").Value = 3336.4446
$ws.Range("L134").Value = 15493.5
$ws.Range("M134").Value = -801.4446000000003
$ws.Range("N134").Value = -20563.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 745.2174
$ws.Range("I131").Value = 390.7143
$ws.Range("J131").Value = 900.3125
$ws.Range("K131").Value = 1172.1429
$ws.Range("L131").Value = 2700.9375
$ws.Range("M131").Value = 3867.8571
$ws.Range("N131").Value = -12780.9375

$ws.Range("H132").Value = 3059
$ws.Range("I132").Value = 869.8
$ws.Range("J132").Value = 4883.3335
$ws.Range("K132").Value = 7828.2
$ws.Range("L132").Value = 43950.0015
$ws.Range("M132").Value = -5298.2
$ws.Range("N132").Value = -49010.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4955.45
$ws.Range("I70").Value = 4914.25
$ws.Range("J70").Value = 4982.9165
$ws.Range("K70").Value = 4914.25
$ws.Range("L70").Value = 4982.9165
$ws.Range("M70").Value = -4644.25
$ws.Range("N70").Value = -5522.9165

$ws.Range("H73").Value = 4955.45
$ws.Range("I73").Value = 4914.25
$ws.Range("J73").Value = 4982.9165
$ws.Range("K73").Value = 4914.25
$ws.Range("L73").Value = 4982.9165
$ws.Range("M73").Value = -3978.25
$ws.Range("N73").Value = -6854.9165

$ws.Range("H102").Value = 2485.4167
$ws.Range("I102").Value = 2099.875
$ws.Range("K102").Value = 2099.875
$ws.Range("M102").Value = -477.875

$ws.Range("H126").Value = 112374.664
$ws.Range("I126").Value = 500746
$ws.Range("J126").Value = 1411.4286
$ws.Range("K126").Value = 1502238
$ws.Range("L126").Value = 4234.2858
$ws.Range("M126").Value = -1499768
$ws.Range("N126").Value = -9174.2858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2537.0833
$ws.Range("I7").Value = 1910
$ws.Range("J7").Value = 3164.1667
$ws.Range("K7").Value = 1910
$ws.Range("L7").Value = 3164.1667
$ws.Range("M7").Value = -1798
$ws.Range("N7").Value = -3388.1667

$ws.Range("H55").Value = 499.84616
$ws.Range("I55").Value = 528.2857
$ws.Range("J55").Value = 466.66666
$ws.Range("K55").Value = 528.2857
$ws.Range("L55").Value = 466.66666
$ws.Range("M55").Value = -355.2857
$ws.Range("N55").Value = -812.66666

$ws.Range("H122").Value = 2826.1667
$ws.Range("I122").Value = 2667.1428
$ws.Range("J122").Value = 3197.2222
$ws.Range("K122").Value = 8001.428400000001
$ws.Range("L122").Value = 9591.6666
$ws.Range("M122").Value = -5551.428400000001
$ws.Range("N122").Value = -14491.6666

$ws.Range("H126").Value = 2537.0833
$ws.Range("I126").Value = 1910
$ws.Range("J126").Value = 3164.1667
$ws.Range("K126").Value = 5730
$ws.Range("L126").Value = 9492.500100000001
$ws.Range("M126").Value = -3260
$ws.Range("N126").Value = -14432.5001

$ws.Range("H133").Value = 28236.223
$ws.Range("J133").Value = 28236.223
$ws.Range("L133").Value = 28236.223
$ws.Range("N133").Value = -33296.223

$ws.Range("H134").Value = 14985.571
$ws.Range("I134").Value = 5000
$ws.Range("J134").Value = 18979.8
$ws.Range("K134").Value = 5000
$ws.Range("L134").Value = 18979.8
$ws.Range("M134").Value = 70
$ws.Range("N134").Value = -29119.8

$ws.Range("H135").Value = 37360
$ws.Range("J135").Value = 37360
$ws.Range("L135").Value = 37360
$ws.Range("N135").Value = -47500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2113.842
$ws.Range("I132").Value = 1050.5555
$ws.Range("J132").Value = 3070.8
$ws.Range("K132").Value = 3151.6665
$ws.Range("L132").Value = 9212.400000000001
$ws.Range("M132").Value = -621.6664999999998
$ws.Range("N132").Value = -14272.4

$ws.Range("H136").Value = 440.54166
$ws.Range("I136").Value = 308.65
$ws.Range("J136").Value = 1100
$ws.Range("K136").Value = 925.9499999999999
$ws.Range("L136").Value = 3300
$ws.Range("M136").Value = 1624.05
$ws.Range("N136").Value = -8400
